# Updated BGR model - 2025-07-29 10:14
#
# The "Electricity Trade Data (TWh) - Source: UNSD" block (header + ISO/attribute
# header + BGR Export row + BGR Import row), previously sitting at rows 80-83
# with a trailing blank spacer row at 84, needs to move down by 5 rows so the
# block now occupies rows 85-88, with the blank spacer row ending up just above
# it at row 84 (and nothing left dangling below row 88).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("historical_data")

# Insert 5 blank rows above the block (80:84); this pushes the existing block
# (header/ISO-header/Export/Import/spacer) down to occupy rows 85:89.
$ws.Rows("80:84").Insert()

# The old trailing spacer row (originally row 84) is now at row 89 - remove it,
# since in the target layout there is no blank row after the Import row (88).
$ws.Rows("89").Delete()

# Give the now-empty row 84 (the row immediately above the shifted block) the
# same 15.75pt custom row height the spacer row originally had.
$ws.Rows("84").RowHeight = 15.75

# Refresh the sheet's used-range selection/dimension to reflect the new extent.
$ws.Activate()
$ws.Range("A1:Z88").Select()
